$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows continuing the date / hours table (rows 8-11)
$ws.Cells.Item(8, 1).Value = 41554
$ws.Cells.Item(8, 2).Value = 4

$ws.Cells.Item(9, 1).Value = 41555
$ws.Cells.Item(9, 2).Value = 2

$ws.Cells.Item(10, 1).Value = 41556
$ws.Cells.Item(10, 2).Value = 6

$ws.Cells.Item(11, 1).Value = 41557
$ws.Cells.Item(11, 2).Value = 1

# Copy the date number-format style from the existing last date cell (A7)
# onto the newly added date cells, so they reuse the same style index
# instead of creating duplicate styles.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8:A11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the active selection, matching the authored workbook state
$ws.Range("A12").Select()
